$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 13.8534432268976
$ws.Range("C2").Value = 10.39168442072861
$ws.Range("D2").Value = 5.976770098523475
$ws.Range("E2").Value = 10.52790139354394
$ws.Range("G2").Value = 35.46075646796771
$ws.Range("H2").Value = 15.12824060728165
$ws.Range("I2").Value = 21.07315753841517
$ws.Range("M2").Value = 15.55281626978474
$ws.Range("N2").Value = 17.2151699519957
$ws.Range("B3").Value = 13.25581190869262
$ws.Range("C3").Value = 9.772851118372374
$ws.Range("D3").Value = 5.857385386439433
$ws.Range("E3").Value = 10.43374717795808
$ws.Range("G3").Value = 34.85203016542171
$ws.Range("H3").Value = 15.11260709436559
$ws.Range("I3").Value = 21.07557214399153
$ws.Range("M3").Value = 15.26405463939967
$ws.Range("N3").Value = 17.28817544722711
$ws.Range("B4").Value = 12.87808856447221
$ws.Range("C4").Value = 9.374074302528651
$ws.Range("D4").Value = 5.78477962777772
$ws.Range("E4").Value = 10.37912318997537
$ws.Range("G4").Value = 34.48815133607764
$ws.Range("H4").Value = 15.10718613313491
$ws.Range("I4").Value = 21.08341602865394
$ws.Range("M4").Value = 15.08857122033334
$ws.Range("N4").Value = 17.33490896383487
$ws.Range("B5").Value = 12.72169906476247
$ws.Range("C5").Value = 9.206970095734441
$ws.Range("D5").Value = 5.755412269310106
$ws.Range("E5").Value = 10.35768466770281
$ws.Range("G5").Value = 34.34257007661797
$ws.Range("H5").Value = 15.10602661655719
$ws.Range("I5").Value = 21.08820567436098
$ws.Range("M5").Value = 15.01761936790536
$ws.Range("N5").Value = 17.35443467458533
$ws.Range("B6").Value = 12.69558965139847
$ws.Range("C6").Value = 9.17894917364336
$ws.Range("D6").Value = 5.75055044827999
$ws.Range("E6").Value = 10.35417495593601
$ws.Range("G6").Value = 34.31856586641615
$ws.Range("H6").Value = 15.10589742310776
$ws.Range("I6").Value = 21.08909700751125
$ws.Range("M6").Value = 15.00587461711322
$ws.Range("N6").Value = 17.35770603004626
$ws.Range("B7").Value = 12.87598907158137
$ws.Range("C7").Value = 9.371839105521259
$ws.Range("D7").Value = 5.784382618794993
$ws.Range("E7").Value = 10.37883071287314
$ws.Range("G7").Value = 34.48617675204036
$ws.Range("H7").Value = 15.10716624796613
$ws.Range("I7").Value = 21.08347418211684
$ws.Range("M7").Value = 15.08761194030353
$ws.Range("N7").Value = 17.33517034296781
$ws.Range("B8").Value = 13.64975422698242
$ws.Range("C8").Value = 10.18228467416116
$ws.Range("D8").Value = 5.935488388349103
$ws.Range("E8").Value = 10.49478798266795
$ws.Range("G8").Value = 35.24894893084165
$ws.Range("H8").Value = 15.12198264669935
$ws.Range("I8").Value = 21.07266628777021
$ws.Range("M8").Value = 15.4529351922359
$ws.Range("N8").Value = 17.23994729035712
$ws.Range("B9").Value = 15.07232480010911
$ws.Range("C9").Value = 11.61788279233176
$ws.Range("D9").Value = 6.235392379426212
$ws.Range("E9").Value = 10.74651031730509
$ws.Range("G9").Value = 36.81336013002818
$ws.Range("H9").Value = 15.18419714897792
$ws.Range("I9").Value = 21.10218229397949
$ws.Range("M9").Value = 16.17930408712014
$ws.Range("N9").Value = 17.06827468783006
$ws.Range("B10").Value = 16.04918826986557
$ws.Range("C10").Value = 12.63120796342894
$ws.Range("D10").Value = 6.455383881747771
$ws.Range("E10").Value = 10.94488958093016
$ws.Range("G10").Value = 37.99133818497936
$ws.Range("H10").Value = 15.25007749246841
$ws.Range("I10").Value = 21.15503071520835
$ws.Range("M10").Value = 16.7132102905112
$ws.Range("N10").Value = 16.95121526620711
$ws.Range("B11").Value = 16.47702355237211
$ws.Range("C11").Value = 13.07268705004991
$ws.Range("D11").Value = 6.554880434645476
$ws.Range("E11").Value = 11.03773192306982
$ws.Range("G11").Value = 38.53066958106547
$ws.Range("H11").Value = 15.28439960952713
$ws.Range("I11").Value = 21.18586704030343
$ws.Range("M11").Value = 16.95497053354235
$ws.Range("N11").Value = 16.89990693365225
$ws.Range("B12").Value = 16.63653781546016
$ws.Range("C12").Value = 13.23602098454992
$ws.Range("D12").Value = 6.592434220959102
$ws.Range("E12").Value = 11.07323401385883
$ws.Range("G12").Value = 38.73517389840377
$ws.Range("H12").Value = 15.29801842475585
$ws.Range("I12").Value = 21.19852116586832
$ws.Range("M12").Value = 17.04626509733701
$ws.Range("N12").Value = 16.88075535569419
$ws.Range("B13").Value = 16.60229634563528
$ws.Range("C13").Value = 13.20101477830926
$ws.Range("D13").Value = 6.584352445965884
$ws.Range("E13").Value = 11.06557314012479
$ws.Range("G13").Value = 38.69112170236965
$ws.Range("H13").Value = 15.29505779165783
$ws.Range("I13").Value = 21.19575243287252
$ws.Range("M13").Value = 17.02661596368965
$ws.Range("N13").Value = 16.88486766144284
$ws.Range("B14").Value = 16.49019755864675
$ws.Range("C14").Value = 13.08620157098473
$ws.Range("D14").Value = 6.557972702766116
$ws.Range("E14").Value = 11.04064596108287
$ws.Range("G14").Value = 38.54749011091295
$ws.Range("H14").Value = 15.2855076034775
$ws.Range("I14").Value = 21.18688852584443
$ws.Range("M14").Value = 16.96248702236119
$ws.Range("N14").Value = 16.89832576375115
$ws.Range("B15").Value = 16.42120533581105
$ws.Range("C15").Value = 13.01537506414438
$ws.Range("D15").Value = 6.541797122659175
$ws.Range("E15").Value = 11.02542136506419
$ws.Range("G15").Value = 38.45954033152846
$ws.Range("H15").Value = 15.27973867162444
$ws.Range("I15").Value = 21.1815863465529
$ws.Range("M15").Value = 16.92317029267102
$ws.Range("N15").Value = 16.90660536801752
$ws.Range("B16").Value = 16.02088474902363
$ws.Range("C16").Value = 12.60181720462346
$ws.Range("D16").Value = 6.448866419170423
$ws.Range("E16").Value = 10.93887194552251
$ws.Range("G16").Value = 37.95614174758981
$ws.Range("H16").Value = 15.24792171729283
$ws.Range("I16").Value = 21.15315225108326
$ws.Range("M16").Value = 16.69738006554197
$ws.Range("N16").Value = 16.95460734199475
$ws.Range("B17").Value = 15.77097323406336
$ws.Range("C17").Value = 12.34124420457308
$ws.Range("D17").Value = 6.391678959048126
$ws.Range("E17").Value = 10.88642186073509
$ws.Range("G17").Value = 37.64804644608358
$ws.Range("H17").Value = 15.22951524751659
$ws.Range("I17").Value = 21.13744965696556
$ws.Range("M17").Value = 16.5585106676446
$ws.Range("N17").Value = 16.98455142278258
$ws.Range("B18").Value = 15.62567929722182
$ws.Range("C18").Value = 12.18883902691744
$ws.Range("D18").Value = 6.358734049006445
$ws.Range("E18").Value = 10.85650033850177
$ws.Range("G18").Value = 37.47117732476043
$ws.Range("H18").Value = 15.21933834950213
$ws.Range("I18").Value = 21.12905763545139
$ws.Range("M18").Value = 16.47853695361048
$ws.Range("N18").Value = 17.00195740459522
$ws.Range("B19").Value = 15.57622272108694
$ws.Range("C19").Value = 12.13680204584327
$ws.Range("D19").Value = 6.347571763441273
$ws.Range("E19").Value = 10.84641261243462
$ws.Range("G19").Value = 37.4113577555033
$ws.Range("H19").Value = 15.21596315869172
$ws.Range("I19").Value = 21.12632609028179
$ws.Range("M19").Value = 16.45144504838013
$ws.Range("N19").Value = 17.00788223869607
$ws.Range("B20").Value = 15.79773822690451
$ws.Range("C20").Value = 12.36924439866077
$ws.Range("D20").Value = 6.397772361752796
$ws.Range("E20").Value = 10.89197997022125
$ws.Range("G20").Value = 37.68081028276097
$ws.Range("H20").Value = 15.23143223982001
$ws.Range("I20").Value = 21.13905501554485
$ws.Range("M20").Value = 16.57330456884749
$ws.Range("N20").Value = 16.98134490232736
$ws.Range("B21").Value = 16.52319233094841
$ws.Range("C21").Value = 13.12002917903645
$ws.Range("D21").Value = 6.565724731254445
$ws.Range("E21").Value = 11.04795855643271
$ws.Range("G21").Value = 38.58967260019917
$ws.Range("H21").Value = 15.28829588818515
$ws.Range("I21").Value = 21.18946555970107
$ws.Range("M21").Value = 16.98133089051302
$ws.Range("N21").Value = 16.89436526660898
$ws.Range("B22").Value = 16.98271750441252
$ws.Range("C22").Value = 13.5883101546467
$ws.Range("D22").Value = 6.674754977935797
$ws.Range("E22").Value = 11.15189498791216
$ws.Range("G22").Value = 39.18514798237062
$ws.Range("H22").Value = 15.32908116776201
$ws.Range("I22").Value = 21.22810596335534
$ws.Range("M22").Value = 17.24647412651222
$ws.Range("N22").Value = 16.83913740538651
$ws.Range("B23").Value = 16.73883085283013
$ws.Range("C23").Value = 13.34042295447347
$ws.Range("D23").Value = 6.616643564227753
$ws.Range("E23").Value = 11.0962492992945
$ws.Range("G23").Value = 38.86726893493108
$ws.Range("H23").Value = 15.30698353175561
$ws.Range("I23").Value = 21.20696219801017
$ws.Range("M23").Value = 17.10513172280838
$ws.Range("N23").Value = 16.86846598982228
$ws.Range("B24").Value = 15.78564279942215
$ws.Range("C24").Value = 12.3565936056075
$ws.Range("D24").Value = 6.395017740273204
$ws.Range("E24").Value = 10.88946642182862
$ws.Range("G24").Value = 37.66599692815827
$ws.Range("H24").Value = 15.23056430500863
$ws.Range("I24").Value = 21.13832725303108
$ws.Range("M24").Value = 16.56661666203086
$ws.Range("N24").Value = 16.982793975959
$ws.Range("B25").Value = 14.69881314756188
$ws.Range("C25").Value = 11.24649810382023
$ws.Range("D25").Value = 6.154143285378183
$ws.Range("E25").Value = 10.67594633909622
$ws.Range("G25").Value = 36.38416765993618
$ws.Range("H25").Value = 15.16381730968904
$ws.Range("I25").Value = 21.08873798707014
$ws.Range("M25").Value = 15.98234689642252
$ws.Range("N25").Value = 17.11311608575199
